$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item(2)
$tf = $shp.TextFrame
$tr = $tf.TextRange

# The "Product Evaluations" bullet is the last paragraph in this placeholder.
$para = $tr.Paragraphs($tr.Paragraphs().Count)

# Split "Product Evaluations" into "Product " + "Evalutions" (typo),
# matching two separate runs as in the target deck.
$word = $para.Characters(9, 11)
$word.Text = "Evalutions"
